$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for A1:E16 (First name, Last name, Month, Day, Year)
$data = @(
    @("Dylan",     "Haworth",   11, 14, 2002),
    @("Dylan",     "DeChoch",    9, 18, 2003),
    @("David",     "Guzman",     1,  2, 2003),
    @("Kaydin",    "Brooksher", 12,  9, 2002),
    @("Christian", "",           2,  9, 2004),
    @("Joel",      "Rocha",      3, 21, 2003),
    @("Thomas",    "Menjivar",   5, 31, 2003),
    @("Trevor",    "",          10, 26, 2002),
    @("John",      "Risso",      9,  4, 2004),
    @("Logan",     "Enslow",     9, 11, 2002),
    @("Fats",      "Waller",    12, 15, 1943),
    @("Jack",      "Haworth",   11, 14, 2002),
    @("Aidan",     "Choi",       4, 12, 2003),
    @("Alex",      "Choi",       1, 27, 2003),
    @("Test",      "Joe",        9, 23, 2000),
    @("John",      "Testies",    9, 23, 2001)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]

    if ([string]::IsNullOrEmpty($vals[1])) {
        $ws.Cells.Item($row, 2).Value = $null
    } else {
        $ws.Cells.Item($row, 2).Value = $vals[1]
    }

    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}

$ws.Range("E16").Select()
